$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# --- Insert a new "Week_Start_Date" column after "Week" (column A) ---
# This shifts ASIN, MyForecast, Amazon Mean/P70/P80/P90 Forecast, Product Title,
# is_holiday_week one column to the right (B->C, C->D, ... H->I, I->J).
$ws.Columns.Item(2).Insert()

# New header for the inserted column
$ws.Range("B1").Value = "Week_Start_Date"

# --- Week labels: drop the leading zero (W01 -> W1, ... W16) ---
$weeks = @("W1","W2","W3","W4","W5","W6","W7","W8","W9","W10","W11","W12","W13","W14","W15","W16")
for ($i = 0; $i -lt $weeks.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $weeks[$i]
}

# --- Week start dates (stored as literal text, not Excel dates) ---
$dates = @(
    "2025-01-05","2025-01-12","2025-01-19","2025-01-26",
    "2025-02-02","2025-02-09","2025-02-16","2025-02-23",
    "2025-03-02","2025-03-09","2025-03-16","2025-03-23",
    "2025-03-30","2025-04-06","2025-04-13","2025-04-20"
)
for ($i = 0; $i -lt $dates.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = "'" + $dates[$i]
}

# --- Updated forecast figures (MyForecast / Amazon Mean / P70 / P80 / P90) ---
# columns D..H, rows 2..17
$forecast = @(
    @(45,38,45,52,62),
    @(46,41,49,58,71),
    @(46,41,50,58,71),
    @(46,43,52,61,75),
    @(47,43,52,62,79),
    @(46,42,51,61,76),
    @(48,44,54,66,85),
    @(49,46,56,70,93),
    @(47,43,53,64,82),
    @(48,44,54,67,88),
    @(49,44,53,67,90),
    @(51,44,54,69,93),
    @(52,43,53,67,89),
    @(51,42,52,66,89),
    @(50,42,52,66,91),
    @(50,41,50,65,88)
)
for ($i = 0; $i -lt $forecast.Length; $i++) {
    $row = $i + 2
    $vals = $forecast[$i]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
    $ws.Cells.Item($row, 6).Value = $vals[2]
    $ws.Cells.Item($row, 7).Value = $vals[3]
    $ws.Cells.Item($row, 8).Value = $vals[4]
}

# --- is_holiday_week becomes a real boolean column ---
$ws.Range("J2:J17").Value = $false

# --- Summary sheet forecast totals (kept as text, matching existing column) ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B9").Value = "'772"
$summary.Range("B10").Value = "'374"
$summary.Range("B11").Value = "'183"
